# Generate Report for Archive
# Update status text "Ready for handoff" -> "In Translation" across sheets,
# and shrink the associated status columns to their new (narrower) width.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: status values live in columns E (zh-cn) and F (de-de), row 2
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E:E").ColumnWidth = 13.4101845877511
$wsOverview.Range("F:F").ColumnWidth = 13.4101845877511

# zh-cn sheet: status value lives in column C, row 2
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C:C").ColumnWidth = 13.4101845877511

# de-de sheet: status value lives in column C, row 2
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C:C").ColumnWidth = 13.4101845877511
